# "workers are now called nodes in slide deck"
# Replace the standalone "Worker" label with "Nodes" wherever it appears
# (the three diagram shapes that label the worker-node boxes).

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "Worker") {
                    $tr.Text = "Nodes"
                }
            }
        }
    }
}
